$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H51").Value = 5000.357
$ws_ALC.Range("I51").Value = 1767.1666
$ws_ALC.Range("K51").Value = 1767.1666
$ws_ALC.Range("M51").Value = -1283.1666
$ws_ALC.Range("H93").Value = 29195
$ws_ALC.Range("J93").Value = 29195
$ws_ALC.Range("L93").Value = 29195
$ws_ALC.Range("N93").Value = -34187
$ws_ALC.Range("H101").Value = 1084.2
$ws_ALC.Range("I101").Value = 560.2222
$ws_ALC.Range("J101").Value = 5800
$ws_ALC.Range("K101").Value = 1680.6666
$ws_ALC.Range("L101").Value = 17400
$ws_ALC.Range("M101").Value = -58.66660000000002
$ws_ALC.Range("N101").Value = -20644
$ws_ALC.Range("H125").Value = 1075.6428
$ws_ALC.Range("I125").Value = 897.4
$ws_ALC.Range("J125").Value = 1174.6666
$ws_ALC.Range("K125").Value = 8076.599999999999
$ws_ALC.Range("L125").Value = 10571.9994
$ws_ALC.Range("M125").Value = -5616.599999999999
$ws_ALC.Range("N125").Value = -15491.9994
$ws_ALC.Range("H129").Value = 965.8
$ws_ALC.Range("I129").Value = 536.1667
$ws_ALC.Range("K129").Value = 1608.5001
$ws_ALC.Range("M129").Value = 3391.4999

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H22").Value = 4081.3333
$ws_ARM.Range("I22").Value = 2275.6
$ws_ARM.Range("J22").Value = 13110
$ws_ARM.Range("K22").Value = 2275.6
$ws_ARM.Range("L22").Value = 13110
$ws_ARM.Range("M22").Value = -1976.6
$ws_ARM.Range("N22").Value = -13708

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 878823.4
$ws_CUL.Range("I2").Value = 36.666668
$ws_CUL.Range("J2").Value = 1255446.2
$ws_CUL.Range("K2").Value = 220.000008
$ws_CUL.Range("L2").Value = 7532677.199999999
$ws_CUL.Range("M2").Value = -107.000008
$ws_CUL.Range("N2").Value = -7532903.199999999
$ws_CUL.Range("H9").Value = 1585.5714
$ws_CUL.Range("I9").Value = 1349.5
$ws_CUL.Range("J9").Value = 1680
$ws_CUL.Range("K9").Value = 4048.5
$ws_CUL.Range("L9").Value = 5040
$ws_CUL.Range("M9").Value = -3824.5
$ws_CUL.Range("N9").Value = -5488
$ws_CUL.Range("H15").Value = 475.29413
$ws_CUL.Range("I15").Value = 160.57143
$ws_CUL.Range("J15").Value = 695.6
$ws_CUL.Range("K15").Value = 481.71429
$ws_CUL.Range("L15").Value = 2086.8
$ws_CUL.Range("M15").Value = -341.71429
$ws_CUL.Range("N15").Value = -2366.8
$ws_CUL.Range("H20").Value = 3340
$ws_CUL.Range("J20").Value = 4933.3335
$ws_CUL.Range("L20").Value = 14800.0005
$ws_CUL.Range("N20").Value = -15254.0005
$ws_CUL.Range("H21").Value = 3670.44
$ws_CUL.Range("I21").Value = 2635.7778
$ws_CUL.Range("J21").Value = 3897.561
$ws_CUL.Range("K21").Value = 7907.3334
$ws_CUL.Range("L21").Value = 11692.683
$ws_CUL.Range("M21").Value = -7734.3334
$ws_CUL.Range("N21").Value = -12038.683
$ws_CUL.Range("H80").Value = 5500.3335
$ws_CUL.Range("J80").Value = 5500.3335
$ws_CUL.Range("L80").Value = 16501.0005
$ws_CUL.Range("N80").Value = -18373.0005
$ws_CUL.Range("H81").Value = 3570
$ws_CUL.Range("I81").Value = 2000
$ws_CUL.Range("J81").Value = 3794.2856
$ws_CUL.Range("K81").Value = 6000
$ws_CUL.Range("L81").Value = 11382.8568
$ws_CUL.Range("M81").Value = -4877
$ws_CUL.Range("N81").Value = -13628.8568
$ws_CUL.Range("H82").Value = 2418.8333
$ws_CUL.Range("I82").Value = 1337.6666
$ws_CUL.Range("J82").Value = 3500
$ws_CUL.Range("K82").Value = 4012.9998
$ws_CUL.Range("L82").Value = 10500
$ws_CUL.Range("M82").Value = -3606.9998
$ws_CUL.Range("N82").Value = -11312
$ws_CUL.Range("H83").Value = 5500.3335
$ws_CUL.Range("J83").Value = 5500.3335
$ws_CUL.Range("L83").Value = 49503.0015
$ws_CUL.Range("N83").Value = -58863.0015
$ws_CUL.Range("H84").Value = 3570
$ws_CUL.Range("I84").Value = 2000
$ws_CUL.Range("J84").Value = 3794.2856
$ws_CUL.Range("K84").Value = 18000
$ws_CUL.Range("L84").Value = 34148.5704
$ws_CUL.Range("M84").Value = -12384
$ws_CUL.Range("N84").Value = -45380.5704
$ws_CUL.Range("H85").Value = 2418.8333
$ws_CUL.Range("I85").Value = 1337.6666
$ws_CUL.Range("J85").Value = 3500
$ws_CUL.Range("K85").Value = 4012.9998
$ws_CUL.Range("L85").Value = 10500
$ws_CUL.Range("M85").Value = -2608.9998
$ws_CUL.Range("N85").Value = -13308
$ws_CUL.Range("H93").Value = 4015.875
$ws_CUL.Range("I93").Value = 0
$ws_CUL.Range("J93").Value = 4015.875
$ws_CUL.Range("K93").Value = 0
$ws_CUL.Range("L93").Value = 12047.625
$ws_CUL.Range("M93").ClearContents()
$ws_CUL.Range("N93").Value = -15791.625
$ws_CUL.Range("H94").Value = 5250
$ws_CUL.Range("J94").Value = 6642.857
$ws_CUL.Range("L94").Value = 19928.571
$ws_CUL.Range("N94").Value = -21280.571
$ws_CUL.Range("H131").Value = 889.14
$ws_CUL.Range("J131").Value = 895.04083
$ws_CUL.Range("L131").Value = 2685.12249
$ws_CUL.Range("N131").Value = -12765.12249
$ws_CUL.Range("H133").Value = 4684.375
$ws_CUL.Range("I133").Value = 3621.25
$ws_CUL.Range("J133").Value = 10000
$ws_CUL.Range("K133").Value = 10863.75
$ws_CUL.Range("L133").Value = 30000
$ws_CUL.Range("M133").Value = -5803.75
$ws_CUL.Range("N133").Value = -40120
$ws_CUL.Range("H134").Value = 5324.645
$ws_CUL.Range("I134").Value = 3196.1667
$ws_CUL.Range("J134").Value = 8271.77
$ws_CUL.Range("K134").Value = 9588.500100000001
$ws_CUL.Range("L134").Value = 24815.31
$ws_CUL.Range("M134").Value = -4518.500100000001
$ws_CUL.Range("N134").Value = -34955.31
$ws_CUL.Range("H137").Value = 20508144
$ws_CUL.Range("I137").Value = 1763.9333
$ws_CUL.Range("J137").Value = 37596790
$ws_CUL.Range("K137").Value = 5291.7999
$ws_CUL.Range("L137").Value = 112790370
$ws_CUL.Range("M137").Value = -191.7999
$ws_CUL.Range("N137").Value = -112800570
$ws_CUL.Range("H139").Value = 1654.2667
$ws_CUL.Range("I139").Value = 1165.1818
$ws_CUL.Range("J139").Value = 2999.25
$ws_CUL.Range("K139").Value = 3495.5454
$ws_CUL.Range("L139").Value = 8997.75
$ws_CUL.Range("M139").Value = 1644.4546
$ws_CUL.Range("N139").Value = -19277.75

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 3226.4666
$ws_GSM.Range("I80").Value = 2610.889
$ws_GSM.Range("J80").Value = 4149.8335
$ws_GSM.Range("K80").Value = 2610.889
$ws_GSM.Range("L80").Value = 4149.8335
$ws_GSM.Range("M80").Value = -1612.889
$ws_GSM.Range("N80").Value = -6145.8335
$ws_GSM.Range("H83").Value = 3226.4666
$ws_GSM.Range("I83").Value = 2610.889
$ws_GSM.Range("J83").Value = 4149.8335
$ws_GSM.Range("K83").Value = 13054.445
$ws_GSM.Range("L83").Value = 20749.1675
$ws_GSM.Range("M83").Value = -8062.445
$ws_GSM.Range("N83").Value = -30733.1675
$ws_GSM.Range("H132").Value = 3090.6428
$ws_GSM.Range("I132").Value = 2789.125
$ws_GSM.Range("J132").Value = 4899.75
$ws_GSM.Range("K132").Value = 8367.375
$ws_GSM.Range("L132").Value = 14699.25
$ws_GSM.Range("M132").Value = -5837.375
$ws_GSM.Range("N132").Value = -19759.25
$ws_GSM.Range("H138").Value = 0
$ws_GSM.Range("J138").Value = 0
$ws_GSM.Range("L138").Value = 0
$ws_GSM.Range("N138").ClearContents()
$ws_GSM.Range("H139").Value = 42346.855
$ws_GSM.Range("J139").Value = 42346.855
$ws_GSM.Range("L139").Value = 42346.855
$ws_GSM.Range("N139").Value = -52626.855

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 628.1875
$ws_LTW.Range("I16").Value = 628.1875
$ws_LTW.Range("J16").Value = 0
$ws_LTW.Range("K16").Value = 628.1875
$ws_LTW.Range("L16").Value = 0
$ws_LTW.Range("M16").Value = -458.1875
$ws_LTW.Range("N16").ClearContents()
$ws_LTW.Range("H18").Value = 0
$ws_LTW.Range("J18").Value = 0
$ws_LTW.Range("L18").Value = 0
$ws_LTW.Range("N18").ClearContents()

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H3").Value = 180000
$ws_WVR.Range("I3").Value = 20000
$ws_WVR.Range("J3").Value = 212000
$ws_WVR.Range("K3").Value = 20000
$ws_WVR.Range("L3").Value = 212000
$ws_WVR.Range("M3").Value = -19886
$ws_WVR.Range("N3").Value = -212228
$ws_WVR.Range("H62").Value = 3660
$ws_WVR.Range("J62").Value = 2980
$ws_WVR.Range("L62").Value = 2980
$ws_WVR.Range("N62").Value = -4228
$ws_WVR.Range("H65").Value = 3660
$ws_WVR.Range("J65").Value = 2980
$ws_WVR.Range("L65").Value = 14900
$ws_WVR.Range("N65").Value = -21140
$ws_WVR.Range("H68").Value = 30000
$ws_WVR.Range("J68").Value = 30000
$ws_WVR.Range("L68").Value = 30000
$ws_WVR.Range("N68").Value = -31622
$ws_WVR.Range("H71").Value = 30000
$ws_WVR.Range("J71").Value = 30000
$ws_WVR.Range("L71").Value = 90000
$ws_WVR.Range("N71").Value = -98112
$ws_WVR.Range("H122").Value = 6791291.5
$ws_WVR.Range("I122").Value = 10418741
$ws_WVR.Range("J122").Value = 94462.53999999999
$ws_WVR.Range("K122").Value = 31256223
$ws_WVR.Range("L122").Value = 283387.62
$ws_WVR.Range("M122").Value = -31253773
$ws_WVR.Range("N122").Value = -288287.62
